$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 165 - shifts the old (blank) row 165 and the summary
# rows below it (166-168) down by one, to row 166-169.
$ws.Rows.Item(165).Insert()

# The last existing entry (row 164) now ends earlier (12:30 instead of
# 12:00), so its duration changes too - formulas recalc automatically.
$ws.Range("E164").Value = 0.52083333333333337

# Populate the newly inserted row 165 with the new time entry
# (2014-07-29, 15:00 - 16:00).
$ws.Range("A165").Value = 2014
$ws.Range("B165").Value = 7
$ws.Range("C165").Value = 29
$ws.Range("D165").Value = 0.625
$ws.Range("E165").Value = 0.66666666666666663
$ws.Range("F165").Formula = "=(E165-D165)*24*60"
$ws.Range("G165").Formula = "=F165/60"

# The "sum [min]" total (now on row 167) needs its range extended to
# include the newly added row.
$ws.Range("F167").Formula = "=SUM(F2:F165)"

# Keep the selection in sync with the shifted active cell.
$ws.Range("F165").Select() | Out-Null
